$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 moves from "Handed back: in sync with en-US" to "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 04:46:45"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a507727464bb3a7d47fca5a4c4ed716ddb84feea/e2e/aa8163b8-00cf-449c-b366-9049651d4851.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03edc075c4c5e7d4481300d498285b3db26c2a3b/e2e/aa8163b8-00cf-449c-b366-9049651d4851.md."

# Target stored column width is exactly 40 "characters"; the host's
# ColumnWidth -> stored-width conversion (Excel's usual MDW-based
# quantization) needs an input of 39 + 1/6 to round-trip to a stored 40.
$colWidth = 39 + 1/6

# --- zh-cn sheet: row 3 status/handoff-date/error-detail update, column P widened ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-25 04:46:40"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $colWidth

# --- de-de sheet: row 3 status/handoff-date/error-detail update, column P widened ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-25 04:46:45"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $colWidth
